$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new user row (row 3): a new client "rochi" with password "1"
$ws.Range("A3").Value = "rochi"
$ws.Range("B3").Value = "'1"
$ws.Range("C3").Value = "Cliente"

# Writing "'1" flags the cell as quote-prefixed text so "1" is kept as a
# string (matching the sibling password cells) instead of becoming numeric.
# Reset the style back to the sheet default so no extra quote-prefix
# formatting is left behind on the cell.
$ws.Range("B3").Style = "Normal"

# Update the selected cell to match the workbook's last-saved selection
$ws.Range("B13").Select()
